# Bump the "Förändrad" (Changed) date column (C) by one day for all data
# rows (2-97): serial 45174 (2023-09-05) -> 45175 (2023-09-06).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:C97")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45174) {
        $cell.Value2 = 45175
    }
}
